# goto pass input to mainframe; prtint_table; log when validation only
#
# Update the "test_jump" sheet (sheet6.xml): rework the jump-logic test
# rows to exercise --jumpto with a key argument, add two new rows for the
# key-based jump scenarios, and tidy up the old validate_method /
# validate_logic / validate_key columns (G:I) that are no longer used in
# favor of a single "logic" column (E) and an optional "key" column (F).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("test_jump")

# --- Row 4: wrong path, should NOT jump -------------------------------------
# drop the old validate_method/validate_logic/validate_key columns (G:I) and
# drive the jump straight off the new "logic" column (E) instead.
$ws.Range("E4").Value = "--jumpto(Yes, 5)"

# --- Row 5: none path, will jump --------------------------------------------
$ws.Range("E5").Value = "--jumpto(No, 0)"

# --- fix up the path descriptions now that the logic changed ---------------
$ws.Range("C4").Value = "wrong path, not jump"
$ws.Range("C6").Value = "none path, will jump"

# --- Row 6: none path, will jump --------------------------------------------
$ws.Range("E6").Value = "--jumpto(No, 3)"

# --- Row 7: key-based jump, yes_key -----------------------------------------
$ws.Range("E7").Value = "--jumpto(Key, 4)"
$ws.Range("F7").Value = "yes_key"

# --- Row 8: none path, key is yes, will not jump ----------------------------
$ws.Range("E8").Value = "--jumpto(Key, 4)"
$ws.Range("F8").Value = "yes_key"

# --- Row 9 (new): key-based jump, no_key ------------------------------------
$ws.Range("A9").Value = "7"
$ws.Range("B9").Value = "css"
$ws.Range("C9").Value = "body > div.container-fluid > div > main > div:nth-child(36) > a.btn.btn-primary.btn-lg.active"
$ws.Range("D9").Value = "checkout"
$ws.Range("E9").Value = "--jumpto(Key, 4)"
$ws.Range("F9").Value = "no_key"
$ws.Range("J9").Value = "fail"

# --- Row 10 (new): none path, key is no, will jump --------------------------
$ws.Range("A10").Value = "8"
$ws.Range("B10").Value = "css"
$ws.Range("D10").Value = "checkout"
$ws.Range("E10").Value = "--jumpto(Key, 4)"
$ws.Range("F10").Value = "no_key"
$ws.Range("J10").Value = "fail"

# --- log the validation-only rows, then the very last new path text --------
$ws.Range("C8").Value = "none path, key is yes, will not jump"
$ws.Range("C10").Value = "none path, key is no, will jump"

# --- now that the new "logic"/"key" columns replace G:I, clear the old ones
$ws.Range("G4:I8").Clear()

# --- cosmetic: the "logic" column needs to re-fit now it holds the longer
# "--jumpto(Key, 4)" / "--jumpto(Yes, 5)" text, and the saved selection moves
# to where the author last clicked ------------------------------------------
$ws.Columns("E").ColumnWidth = 12.91
$ws.Range("H8").Select()
